# TC19 Updated for ECTEST
# Set Runmode (column E) to "No" for all test cases except TC19 (row 19),
# which keeps its Runmode value of "Yes".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

$lastRow = 29
for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -ne 19) {
        $ws.Cells.Item($r, 5).Value = "No"
    }
}

# Update the active sheet selection / scroll position to match the
# author's last saved view.
$ws.Range("E24").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
